# Updated symbol list (price/volume refresh) per commit on 2023-01-14
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "304.91"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "6.07%"
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.90"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "7.89%"
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.280"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.93%"
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07327"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "9.22%"
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.822"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "6.53%"
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.769"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "10.97%"
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.447"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "5.82%"
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9094"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.19%"
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01636"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2,426.83%"
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1680"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "5.61%"
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07496"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "9.77%"
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08005"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3.90%"
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03032"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "3.44%"
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09976"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "11.11%"
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001510"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-4.79%"
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04571"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.66%"
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006235"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.25%"
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.491"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.34%"
# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.20%"
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.3322"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "3.31%"
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1324"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.60%"
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.300"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "5.55%"
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1638"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.60%"
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001233"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "3.53%"
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004418"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "7.07%"
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001316"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "9.86%"
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001761"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "9.03%"
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04492"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "5.20%"
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007122"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.74%"
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1350"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "8.77%"
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002288"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "3.72%"
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01428"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "17.96%"
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006084"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "7.68%"
# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.83%"
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01315"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.79%"
